$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Status updates (column F) ---------------------------------------
# Row 6 ("Decide Game rules...") was "In Progress" and is now "Completed".
$ws.Range("F6").Value = "Completed"

# Newly tracked requirements that are marked "Completed".
$completedCells = @("F10", "F11", "F12", "F17", "F19", "F22", "F24", "F25")
foreach ($addr in $completedCells) {
    $ws.Range($addr).Value = "Completed"
}

# Newly tracked requirements that are marked "In Progress" (shown in red).
$inProgressCells = @("F8", "F13", "F18", "F20")
foreach ($addr in $inProgressCells) {
    $ws.Range($addr).Value = "In Progress"
    $ws.Range($addr).Font.Color = 255
}

# --- Selection / view state ------------------------------------------
$ws.Range("E14").Select()
